$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: The two highlighted "question" paragraphs under part (c)
#   "Calculate the heat sink size for a short-circuit output."
#   "Calculate short-circuit current."
# are replaced by four plain (un-highlighted) answer paragraphs.
# ---------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute("Calculate the heat sink size for a short-circuit output.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$qPara1 = $find1.Paragraphs(1)

$find2 = $d.Content
$find2.Find.Execute("Calculate short-circuit current.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$qPara2 = $find2.Paragraphs(1)

$insertAt = $qPara1.Range.Start
$toRemove = $d.Range($qPara1.Range.Start, $qPara2.Range.End)
$toRemove.Delete()

$answerText = "Power dissipation = 10.871 V * 1.6877 A =  18.3 W (short circuit)`r" + `
    "Max Case Temp = ~130 °C`r" + `
    "Heat sink size = (130 °C – 50 °C)/18.3 W = 4.37 C°/W`r" + `
    "Short Circuit Current = 1.678 A`r"

$dest = $d.Range($insertAt, $insertAt)
$dest.InsertAfter($answerText)

# Match the paragraph-mark formatting of the first new paragraph (single
# line spacing) recorded by the author's save.
$firstNewPara = $d.Range($insertAt, $insertAt).Paragraphs(1)
$firstNewPara.Range.ParagraphFormat.LineSpacingRule = 0

# ---------------------------------------------------------------------
# Edit 2: The highlighted prompt sentence under part (e)
#   "Modify resistor values for higher q current and measure load
#    regulation at Vin = 15V"
# is replaced by the resistor-substitution measurement write-up.
# ---------------------------------------------------------------------
$find3 = $d.Content
$find3.Find.Execute("Modify resistor values for higher q current and measure load regulation at Vin = 15V", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$promptPara = $find3.Paragraphs(1)

$insertAt2 = $promptPara.Range.Start
$removeEnd2 = $promptPara.Range.End - 1
$toRemove2 = $d.Range($insertAt2, $removeEnd2)
$toRemove2.Delete()

$measurementText = "Decreasing R3 from 68 Ω to 43 Ω DID improve load regulation.  " + `
    "At no load, output voltage was 8.51 V, but dropped to 8.50 V with a 473 mA load, " + `
    "and to 8.45 V with a 1.0 A load.  With the 43 Ω resistance instead, no load output " + `
    "voltage was 8.6 V, with a 488 mA load, it was 8.57 V, and with a 1.0 A load, 8.54 V.  "

$dest2 = $d.Range($insertAt2, $insertAt2)
$dest2.InsertAfter($measurementText)

Write-Host "Heat sink + resistor substitution answers written."
